$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44428
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 1550
$ws.Range("D3").Value = 44218
$ws.Range("J3").Value = 320
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 10500
$ws.Range("P3").Value = 1050
$ws.Range("D4").Value = 44580
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19000
$ws.Range("P4").Value = 1900
$ws.Range("D5").Value = 44377
$ws.Range("J5").Value = 650
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14538
$ws.Range("P5").Value = 1454
$ws.Range("D6").Value = 44441
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 1550
$ws.Range("D7").Value = 44406
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("P7").Value = 1450
$ws.Range("D8").Value = 44263
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 1550
$ws.Range("D9").Value = 44265
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 1550
$ws.Range("D10").Value = 44204
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 10500
$ws.Range("P10").Value = 1050
$ws.Range("D11").Value = 44358
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 1450
$ws.Range("D12").Value = 44460
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 1550
$ws.Range("D13").Value = 44547
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19500
$ws.Range("P13").Value = 1950
$ws.Range("D14").Value = 44291
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("P14").Value = 1350
$ws.Range("D15").Value = 44524
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 21000
$ws.Range("M15").Value = 20500
$ws.Range("P15").Value = 2050
$ws.Range("D16").Value = 44330
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 13500
$ws.Range("P16").Value = 1350
$ws.Range("D17").Value = 44160
$ws.Range("J17").Value = 360
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 11000
$ws.Range("M17").Value = 10500
$ws.Range("P17").Value = 1050
